# opencloning_linkml.xlsx edit:
#   - add a new "OpenDNACollectionsSource" sheet (copy of the IGEMSource
#     layout) right after "IGEMSource"
#   - append ",open_dna_collections" to the repository-name dropdown list
#     used on every *Source sheet that offers it

$wb = $excel.ActiveWorkbook

$newFormula = '"addgene,genbank,benchling,snapgene,euroscarf,igem,wekwikgene,seva,open_dna_collections"'

# Sheets whose repository-name dropdown needs the new option, and the
# column range the validation lives on (mirrors the source workbook).
$targets = @(
    @{ Sheet = "RepositoryIdSource";     Range = "B2:B1048576" },
    @{ Sheet = "AddgeneIdSource";        Range = "D2:D1048576" },
    @{ Sheet = "WekWikGeneIdSource";     Range = "C2:C1048576" },
    @{ Sheet = "SEVASource";             Range = "C2:C1048576" },
    @{ Sheet = "BenchlingUrlSource";     Range = "B2:B1048576" },
    @{ Sheet = "SnapGenePlasmidSource";  Range = "B2:B1048576" },
    @{ Sheet = "EuroscarfSource";        Range = "B2:B1048576" },
    @{ Sheet = "IGEMSource";             Range = "C2:C1048576" }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $rng = $ws.Range($t.Range)
    $rng.Validation.Delete()
    $rng.Validation.Add(3, 1, 1, $newFormula)
}

# Insert the new OpenDNACollectionsSource sheet right after IGEMSource,
# with the same column layout as IGEMSource (sequence_file_url /
# repository_id / repository_name / type / output_name / database_id /
# input / id) and the same (updated) dropdown on repository_name (col C).
$after = $wb.Worksheets.Item("IGEMSource")
$newSheet = $wb.Worksheets.Add($null, $after)
$newSheet.Name = "OpenDNACollectionsSource"

$headers = @("sequence_file_url", "repository_id", "repository_name", "type", "output_name", "database_id", "input", "id")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$newRange = $newSheet.Range("C2:C1048576")
$newRange.Validation.Add(3, 1, 1, $newFormula)
